$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 266
$endRow = 385

# 1) Snapshot-read the columns that vary per data row (D, J, K, L, M, P)
#    for every row in the existing data block BEFORE any writes happen.
$snapD  = $ws.Range("D$($startRow):D$($endRow)").Value2
$snapJM = $ws.Range("J$($startRow):M$($endRow)").Value2
$snapP  = $ws.Range("P$($startRow):P$($endRow)").Value2

# 2) The two brand new rows at the bottom (386, 387) are exact duplicates of
#    the last two existing rows (384, 385) -- copy every column so the new
#    rows inherit A,B,C,E,F,G,H,I,N,O,Q,R (which are constant/alternating)
#    as well as their own D,J,K,L,M,P values.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(386, $col).Value = $ws.Cells.Item(384, $col).Value2
    $ws.Cells.Item(387, $col).Value = $ws.Cells.Item(385, $col).Value2
}
$ws.Cells.Item(386, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(387, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 3) Shift every existing row down by one pair (two rows): new row r gets
#    the (pre-edit) values that used to live in row r-2. Walk from the
#    bottom upward so we never overwrite a source row before reading it.
for ($r = $endRow; $r -ge ($startRow + 2); $r--) {
    $idx = ($r - 2) - $startRow + 1
    $ws.Cells.Item($r, 4).Value  = $snapD[$idx, 1]
    $ws.Cells.Item($r, 10).Value = $snapJM[$idx, 1]
    $ws.Cells.Item($r, 11).Value = $snapJM[$idx, 2]
    $ws.Cells.Item($r, 12).Value = $snapJM[$idx, 3]
    $ws.Cells.Item($r, 13).Value = $snapJM[$idx, 4]
    $ws.Cells.Item($r, 16).Value = $snapP[$idx, 1]
}

# 4) Insert the new, most-recent pair of readings into rows 266/267.
$ws.Cells.Item(266, 4).Value  = 44845
$ws.Cells.Item(266, 10).Value = 1000

$ws.Cells.Item(267, 4).Value  = 44845
$ws.Cells.Item(267, 10).Value = 1000
$ws.Cells.Item(267, 11).Value = 450
$ws.Cells.Item(267, 12).Value = 500
$ws.Cells.Item(267, 13).Value = 475
$ws.Cells.Item(267, 16).Value = 95
